$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the date for the new row (8): 2021-11-05 -> serial 44505
# Copy style from B7 (already formatted as date with border) so we reuse the
# existing style entry rather than creating a new numFmt/style.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B8").Value = 44505

# Enter the hours worked
$ws.Range("C8").Value = 2

# Move selection to C9, mimicking pressing Enter after data entry
$ws.Range("C9").Select()
